# Add a new "V_section [m³]" column (H) to the Querschnittswerte sheet and
# extend the merged header band (B5:G5 -> B5:H5) to cover the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Formatting for the merged header row (row 5) -----------------------
# G5 currently closes the merged band (border on top/bottom/right).
# Since H becomes the new last column, copy that "closing" format to H5
# first (while G5 still has it), then turn G5 into a regular "middle" cell
# by copying the format already used by C5:F5.
$ws.Range("G5").Copy()
$ws.Range("H5").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("C5").Copy()
$ws.Range("G5").PasteSpecial(-4122)   # xlPasteFormats

# --- 2. Extend the merged cell B5:G5 to B5:H5 -------------------------------
$ws.Range("B5:G5").UnMerge()
$ws.Range("B5:H5").Merge()

# --- 3. New header cell H6 --------------------------------------------------
# Copy the header formatting from the neighboring header cell (G6) and set
# the label text.
$ws.Range("G6").Copy()
$ws.Range("H6").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H6").Value = "V_section [m³]"

# --- 4. New data values in column H (rows 8-17) -----------------------------
$ws.Range("H8").Value = 170.88
$ws.Range("H9").Value = 157.24
$ws.Range("H10").Value = 143.6
$ws.Range("H11").Value = 129.96
$ws.Range("H12").Value = 116.31
$ws.Range("H13").Value = 102.67
$ws.Range("H14").Value = 89.03
$ws.Range("H15").Value = 75.39
$ws.Range("H16").Value = 61.75
$ws.Range("H17").Value = 0
